$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 highlighted (yellow) outlier rows, from bottom to top so
# the remaining row numbers don't shift under us while we work.
$rowsToDelete = @(111, 101, 64, 61, 58)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete() | Out-Null
}

# Restore portrait page orientation.
$ws.PageSetup.Orientation = 1

# Update the view: scrolled down, with row 107 selected (full row).
$ws.Range("A107:XFD107").Select() | Out-Null

Write-Output "done"
